$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2 was blank; give it the same email (shared string) the rest of column B
# already uses, then turn it into a real mailto: hyperlink like B3:B14.
$ws.Range("B2").Value = "aruncyclopse007@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:aruncyclopse007@gmail.com")

# Hyperlinks.Add re-stamps the cell with a fresh (duplicate) style record;
# restore the original "Hyperlink" cell style B3:B14 already use so B2
# keeps sharing that same style index instead of growing the style table.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
